$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(95).Insert()
$ws.Range("B95").Value = "SP(1-0-00)"
$ws.Range("K95").Value = 45120
